$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.165.08"
$ws.Range("E2").Value = "  +4.04%  "
$ws.Range("D3").Value = "3.661.53"
$ws.Range("E3").Value = "  +7.48%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'594.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Value = "'182.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").Value = "3.646.71"
$ws.Range("E7").Value = "  +7.25%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.205"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.87%  "
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'50.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").Value = "'0.0000289"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "'696.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").Value = "4.248.28"
$ws.Range("E15").Value = "  +7.32%  "
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "72.193.11"
$ws.Range("E17").Value = "  +3.84%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.632.97"
$ws.Range("E18").Value = "  +6.53%  "
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("E20").Value = "  +4.75%  "
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("E22").Value = "  +3.09%  "
$ws.Range("D23").Value = "'5.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.11%  "
$ws.Range("D24").Value = "'17.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").Value = "'104.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'4.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("E27").Value = "  +5.17%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("D29").Value = "'35.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.54%  "
$ws.Range("E30").Value = "  +3.79%  "
$ws.Range("D31").Value = "'7.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.88%  "
$ws.Range("E32").Value = "  +16.42%  "
$ws.Range("D33").Value = "'583.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.24%  "
$ws.Range("D34").Value = "'11.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'59.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.37%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "3.673.75"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'0.145"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "'36.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  +7.54%  "
$ws.Range("D42").Value = "'3.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.21%  "
$ws.Range("D43").Value = "'0.0466"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.37%  "
$ws.Range("E44").Value = "  +4.13%  "
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'2.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.11%  "
$ws.Range("E48").Value = "  +2.28%  "
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'132.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
